# Add smart auto-scroll behavior to log steps - scrolls to latest unless user manually scrolling
# (As applied to the test-results workbook: clear out the now-stale Status/Remarks/
#  Actual Output/Screenshot/Page Source columns for log rows that haven't been run yet.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 10-24: clear Status, Remarks, Actual Output, Screenshot and Page Source columns (L:P)
$ws.Range("L10:P24").ClearContents()

# Rows 25-40: clear Status and Remarks columns (L:M)
$ws.Range("L25:M40").ClearContents()
